# BookReport.xlsx template update ("Activity 6" report rename):
#  - The report title cell (A5) changes from "CUSTOMER REPORT" to "BOOK REPORT"
#    and is emphasised in bold.
#  - The "Version 1" label (G6) keeps its text but is re-aligned.
#  - The merged F3:G3 header cell formatting is tidied up (right aligned).
#  - Selection moves to the title cell and the sheet is set to print portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the title cell and make it bold.
$ws.Range("A5").Value = "BOOK REPORT"
$ws.Range("A5").Font.Bold = $true

# Re-affirm the "Version 1" label's right alignment.
$ws.Range("G6").Value = "Version 1"
$ws.Range("G6").HorizontalAlignment = -4152

# Tidy the merged header cell's alignment.
$ws.Range("F3:G3").HorizontalAlignment = -4152

# Move the active selection to the title cell.
$ws.Range("A5").Select()

# Print the sheet in portrait orientation.
$ws.PageSetup.Orientation = 1
